$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.441254854202271
$ws.Range("B1").Value = 3.324631452560425
$ws.Range("C1").Value = 4.370556354522705
$ws.Range("D1").Value = 2.019719362258911
$ws.Range("E1").Value = 1.159057378768921
